$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1216.5
$ws.Range("J112").Value = 1445.625
$ws.Range("L112").Value = 4336.875
$ws.Range("N112").Value = -6552.875
$ws.Range("H137").Value = 1639.8148
$ws.Range("I137").Value = 1215.0526
$ws.Range("K137").Value = 3645.1578
$ws.Range("M137").Value = -1095.1578
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3665
$ws.Range("I3").Value = 3665
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3665
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3550
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 2420.5
$ws.Range("I15").Value = 2223
$ws.Range("J15").Value = 3013
$ws.Range("K15").Value = 2223
$ws.Range("L15").Value = 3013
$ws.Range("M15").Value = -1873
$ws.Range("N15").Value = -3713
$ws.Range("H61").Value = 1878.75
$ws.Range("I61").Value = 1534
$ws.Range("J61").Value = 2223.5
$ws.Range("K61").Value = 1534
$ws.Range("L61").Value = 2223.5
$ws.Range("M61").Value = -1322
$ws.Range("N61").Value = -2647.5
$ws.Range("H74").Value = 9376966
$ws.Range("I74").Value = 13638138
$ws.Range("K74").Value = 13638138
$ws.Range("M74").Value = -13637264
$ws.Range("H77").Value = 9376966
$ws.Range("I77").Value = 13638138
$ws.Range("K77").Value = 68190690
$ws.Range("M77").Value = -68186322
$ws.Range("H88").Value = 142859890
$ws.Range("J88").Value = 250003340
$ws.Range("L88").Value = 250003340
$ws.Range("N88").Value = -250004152
$ws.Range("H91").Value = 142859890
$ws.Range("J91").Value = 250003340
$ws.Range("L91").Value = 250003340
$ws.Range("N91").Value = -250006148
$ws.Range("H132").Value = 4528.095
$ws.Range("I132").Value = 4478.5264
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 13435.5792
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -10905.5792
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 1878.75
$ws.Range("I136").Value = 1534
$ws.Range("J136").Value = 2223.5
$ws.Range("K136").Value = 4602
$ws.Range("L136").Value = 6670.5
$ws.Range("M136").Value = -2052
$ws.Range("N136").Value = -11770.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18183962
$ws.Range("I86").Value = 33334930
$ws.Range("J86").Value = 2800.4
$ws.Range("K86").Value = 33334930
$ws.Range("L86").Value = 2800.4
$ws.Range("M86").Value = -33333807
$ws.Range("N86").Value = -5046.4
$ws.Range("H89").Value = 18183962
$ws.Range("I89").Value = 33334930
$ws.Range("J89").Value = 2800.4
$ws.Range("K89").Value = 166674650
$ws.Range("L89").Value = 14002
$ws.Range("M89").Value = -166669034
$ws.Range("N89").Value = -25234
$ws.Range("H134").Value = 2633.7097
$ws.Range("I134").Value = 2250.6667
$ws.Range("J134").Value = 3164.077
$ws.Range("K134").Value = 6752.000100000001
$ws.Range("L134").Value = 9492.231
$ws.Range("M134").Value = -4217.000100000001
$ws.Range("N134").Value = -14562.231
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8718673
$ws.Range("I31").Value = 6091953
$ws.Range("J31").Value = 15386502
$ws.Range("K31").Value = 6091953
$ws.Range("L31").Value = 15386502
$ws.Range("M31").Value = -6091658
$ws.Range("N31").Value = -15387092
$ws.Range("H34").Value = 8718673
$ws.Range("I34").Value = 6091953
$ws.Range("J34").Value = 15386502
$ws.Range("K34").Value = 6091953
$ws.Range("L34").Value = 15386502
$ws.Range("M34").Value = -6091751
$ws.Range("N34").Value = -15386906
$ws.Range("H58").Value = 1309.425
$ws.Range("I58").Value = 826.56525
$ws.Range("J58").Value = 1962.7059
$ws.Range("K58").Value = 826.56525
$ws.Range("L58").Value = 1962.7059
$ws.Range("M58").Value = -623.56525
$ws.Range("N58").Value = -2368.7059
$ws.Range("H99").Value = 1988464.8
$ws.Range("I99").Value = 3576463
$ws.Range("J99").Value = 3467.125
$ws.Range("K99").Value = 3576463
$ws.Range("L99").Value = 3467.125
$ws.Range("M99").Value = -3574965
$ws.Range("N99").Value = -6463.125
$ws.Range("H126").Value = 1988464.8
$ws.Range("I126").Value = 3576463
$ws.Range("J126").Value = 3467.125
$ws.Range("K126").Value = 10729389
$ws.Range("L126").Value = 10401.375
$ws.Range("M126").Value = -10726919
$ws.Range("N126").Value = -15341.375
$ws.Range("H132").Value = 2001.0435
$ws.Range("I132").Value = 1246.5454
$ws.Range("K132").Value = 3739.6362
$ws.Range("M132").Value = -1209.6362
$ws.Range("H134").Value = 3252.3333
$ws.Range("I134").Value = 3463.4211
$ws.Range("J134").Value = 2751
$ws.Range("K134").Value = 10390.2633
$ws.Range("L134").Value = 8253
$ws.Range("M134").Value = -7855.263300000001
$ws.Range("N134").Value = -13323
$ws.Range("H136").Value = 1309.425
$ws.Range("I136").Value = 826.56525
$ws.Range("J136").Value = 1962.7059
$ws.Range("K136").Value = 2479.69575
$ws.Range("L136").Value = 5888.1177
$ws.Range("M136").Value = 70.30425000000014
$ws.Range("N136").Value = -10988.1177
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H33").Value = 3996.6155
$ws.Range("I33").Value = 252.83333
$ws.Range("J33").Value = 7205.5713
$ws.Range("K33").Value = 1516.99998
$ws.Range("L33").Value = 43233.4278
$ws.Range("M33").Value = -1233.99998
$ws.Range("N33").Value = -43799.4278
$ws.Range("H56").Value = 7110
$ws.Range("I56").Value = 7110
$ws.Range("K56").Value = 7110
$ws.Range("M56").Value = -6580
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 3100
$ws.Range("J24").Value = 3100
$ws.Range("L24").Value = 3100
$ws.Range("N24").Value = -3446
$ws.Range("H132").Value = 3274.0667
$ws.Range("I132").Value = 2391.5557
$ws.Range("J132").Value = 4597.8335
$ws.Range("K132").Value = 7174.6671
$ws.Range("L132").Value = 13793.5005
$ws.Range("M132").Value = -4644.6671
$ws.Range("N132").Value = -18853.5005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17587120
$ws.Range("I132").Value = 33102184
$ws.Range("J132").Value = 3380.0667
$ws.Range("K132").Value = 99306552
$ws.Range("L132").Value = 10140.2001
$ws.Range("M132").Value = -99304022
$ws.Range("N132").Value = -15200.2001
$ws.Range("H136").Value = 5441462.5
$ws.Range("I136").Value = 11375891
$ws.Range("J136").Value = 1570
$ws.Range("K136").Value = 34127673
$ws.Range("L136").Value = 4710
$ws.Range("M136").Value = -34125123
$ws.Range("N136").Value = -9810
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2805.5386
$ws.Range("I132").Value = 1008.625
$ws.Range("J132").Value = 5680.6
$ws.Range("K132").Value = 3025.875
$ws.Range("L132").Value = 17041.8
$ws.Range("M132").Value = -495.875
$ws.Range("N132").Value = -22101.8
$ws.Range("H136").Value = 1266.2084
$ws.Range("I136").Value = 660.5
$ws.Range("J136").Value = 3083.3333
$ws.Range("K136").Value = 1981.5
$ws.Range("L136").Value = 9249.999899999999
$ws.Range("M136").Value = 568.5
$ws.Range("N136").Value = -14349.9999
